$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D that must remain stored as text even though their new value looks numeric
# (NumberFormat forced to Text before assignment to prevent Excel auto-converting to a number)
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'

$ws.Range('D2').Value = '34.641.28'
$ws.Range('E2').Value = '  +1.42%  '
$ws.Range('D3').Value = '1.801.46'
$ws.Range('E3').Value = '  +0.83%  '
$ws.Range('E4').Value = '  -0.33%  '
$ws.Range('D5').Value = '226.87'
$ws.Range('E5').Value = '  +0.14%  '
$ws.Range('D6').Value = '0.560'
$ws.Range('E6').Value = '  +2.28%  '
$ws.Range('D7').Value = '0.997'
$ws.Range('E7').Value = '  -0.37%  '
$ws.Range('E8').Value = '  +3.66%  '
$ws.Range('E9').Value = '  +2.14%  '
$ws.Range('D10').Value = '0.0696'
$ws.Range('E10').Value = '  +1.06%  '
$ws.Range('D11').Value = '0.0949'
$ws.Range('E11').Value = '  +0.38%  '
$ws.Range('D12').Value = '2.058.51'
$ws.Range('E12').Value = '  +0.66%  '
$ws.Range('D13').Value = '11.12'
$ws.Range('E13').Value = '  +0.69%  '
$ws.Range('D14').Value = '1.796.91'
$ws.Range('E14').Value = '  +0.39%  '
$ws.Range('E15').Value = '  +2.51%  '
$ws.Range('D16').Value = '34.590.01'
$ws.Range('D17').Value = '4.29'
$ws.Range('E17').Value = '  +2.58%  '
$ws.Range('D18').Value = '69.04'
$ws.Range('E18').Value = '  +1.12%  '
$ws.Range('D19').Value = '248.44'
$ws.Range('E19').Value = '  +0.82%  '
$ws.Range('D20').Value = '0.0₃0803'
$ws.Range('E20').Value = '  +3.38%  '
$ws.Range('D21').Value = '11.38'
$ws.Range('E21').Value = '  +4.42%  '
$ws.Range('E22').Value = '  -0.33%  '
$ws.Range('E23').Value = '  +1.65%  '
$ws.Range('E24').Value = '  +0.37%  '
$ws.Range('D25').Value = '164.56'
$ws.Range('E25').Value = '  +1.92%  '
$ws.Range('D26').Value = '7.29'
$ws.Range('E26').Value = '  +1.43%  '
$ws.Range('D27').Value = '16.58'
$ws.Range('E27').Value = '  +1.57%  '
$ws.Range('E28').Value = '  +2.71%  '
$ws.Range('E29').Value = '  -0.33%  '
$ws.Range('D30').Value = '3.97'
$ws.Range('E30').Value = '  +9.43%  '
$ws.Range('D31').Value = '3.81'
$ws.Range('E31').Value = '  +3.46%  '
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').Value = '0.0522'
$ws.Range('E32').Value = '  +0.48%  '
$ws.Range('B33').Value = 'PancakeSwap'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D33').Value = '1.23'
$ws.Range('E33').Value = '  -0.42%  '
$ws.Range('E34').Value = '  +1.74%  '
$ws.Range('D35').Value = '1.425.13'
$ws.Range('E35').Value = '  -1.47%  '
$ws.Range('E36').Value = '  +5.55%  '
$ws.Range('D37').Value = '0.672'
$ws.Range('E37').Value = '  +2.74%  '
$ws.Range('E38').Value = '  +0.66%  '
$ws.Range('E39').Value = '  +1.89%  '
$ws.Range('D40').Value = '85.49'
$ws.Range('E40').Value = '  +6.13%  '
$ws.Range('E41').Value = '  +1.87%  '
$ws.Range('E42').Value = '  +0.65%  '
$ws.Range('E43').Value = '  +2.20%  '
$ws.Range('D44').Value = '13.50'
$ws.Range('E44').Value = '  -0.48%  '
$ws.Range('D45').Value = '0.0525'
$ws.Range('E45').Value = '  +3.33%  '
$ws.Range('D46').Value = '6.05'
$ws.Range('E46').Value = '  -0.29%  '
$ws.Range('E47').Value = '  +0.30%  '
$ws.Range('D48').Value = '1.957.10'
$ws.Range('E48').Value = '  +0.53%  '
$ws.Range('D49').Value = '105.87'
$ws.Range('E49').Value = '  -0.02%  '
$ws.Range('E50').Value = '  -0.35%  '
$ws.Range('E51').Value = '  -5.47%  '
